$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Delete entire row 28 (Caso 5738, IRALA 29) - all subsequent rows shift up.
$ws.Rows.Item(28).Delete()
